$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 237
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
